# Update "想去人数" (F column) values across sheets to match regenerated data
# (gh-pages output generated at 456a3b4)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 2763
$ws1.Range("F5").Value  = 1538
$ws1.Range("F6").Value  = 1148
$ws1.Range("F12").Value = 9352
$ws1.Range("F14").Value = 2506
$ws1.Range("F32").Value = 165
$ws1.Range("F41").Value = 1062
$ws1.Range("F43").Value = 1432
$ws1.Range("F45").Value = 317
$ws1.Range("F49").Value = 305

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F10").Value = 2

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 2763
$ws4.Range("F5").Value  = 1538
$ws4.Range("F7").Value  = 1148
$ws4.Range("F10").Value = 9352
$ws4.Range("F12").Value = 2506
$ws4.Range("F26").Value = 165
$ws4.Range("F38").Value = 1062
$ws4.Range("F41").Value = 1432
$ws4.Range("F44").Value = 317
$ws4.Range("F48").Value = 305
$ws4.Range("F49").Value = 2
